# Updates cryptos list price/volume figures (and two coin row swaps) to
# match the latest GitHub Actions scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.517.71"
$ws.Range("E2").Value = "  +1.19%  "

# Row 3
$ws.Range("D3").Value = "3.013.81"
$ws.Range("E3").Value = "  +2.53%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'378.14"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6
$ws.Range("D6").Value = "'103.03"
$ws.Range("E6").Value = "  +2.37%  "

# Row 7
$ws.Range("E7").Value = "  +1.24%  "

# Row 9
$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  +2.51%  "

# Row 10
$ws.Range("D10").Value = "'36.59"
$ws.Range("E10").Value = "  +1.62%  "

# Row 11
$ws.Range("E11").Value = "  -0.35%  "

# Row 12
$ws.Range("D12").Value = "'0.0858"
$ws.Range("E12").Value = "  +0.78%  "

# Row 13
$ws.Range("D13").Value = "3.499.09"
$ws.Range("E13").Value = "  +3.04%  "

# Row 14
$ws.Range("D14").Value = "'18.47"
$ws.Range("E14").Value = "  +1.49%  "

# Row 15
$ws.Range("D15").Value = "'7.72"
$ws.Range("E15").Value = "  +1.47%  "

# Row 16
$ws.Range("D16").Value = "3.012.87"
$ws.Range("E16").Value = "  +2.31%  "

# Row 17
$ws.Range("D17").Value = "'0.978"
$ws.Range("E17").Value = "  -1.63%  "

# Row 18
$ws.Range("D18").Value = "'10.51"
$ws.Range("E18").Value = "  -13.63%  "

# Row 19
$ws.Range("D19").Value = "51.531.31"
$ws.Range("E19").Value = "  +1.28%  "

# Row 20
$ws.Range("D20").Value = "'3.02"
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("D21").Value = "'12.45"
$ws.Range("E21").Value = "  +0.45%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +1.44%  "

# Row 23
$ws.Range("D23").Value = "'69.87"
$ws.Range("E23").Value = "  +0.54%  "

# Row 24
$ws.Range("D24").Value = "'267.09"
$ws.Range("E24").Value = "  +0.12%  "

# Row 25
$ws.Range("D25").Value = "'3.12"
$ws.Range("E25").Value = "  -3.60%  "

# Row 26
$ws.Range("D26").Value = "'8.18"
$ws.Range("E26").Value = "  +2.16%  "

# Row 27
$ws.Range("D27").Value = "'7.50"
$ws.Range("E27").Value = "  +5.66%  "

# Row 28
$ws.Range("E28").Value = "  +5.63%  "

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").Value = "'26.19"
$ws.Range("E30").Value = "  +2.40%  "

# Row 31
$ws.Range("D31").Value = "'0.108"
$ws.Range("E31").Value = "  -0.15%  "

# Row 32
$ws.Range("D32").Value = "'10.28"
$ws.Range("E32").Value = "  +2.62%  "

# Row 33
$ws.Range("D33").Value = "'34.10"
$ws.Range("E33").Value = "  +1.82%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'50.67"
$ws.Range("E34").Value = "  +0.49%  "

# Row 35
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.0455"
$ws.Range("E35").Value = "  +5.28%  "

# Row 36
$ws.Range("E36").Value = "  +0.26%  "

# Row 37
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("E38").Value = "  +6.77%  "

# Row 39
$ws.Range("D39").Value = "'17.25"
$ws.Range("E39").Value = "  +4.27%  "

# Row 40
$ws.Range("D40").Value = "'0.286"
$ws.Range("E40").Value = "  +10.41%  "

# Row 41
$ws.Range("D41").Value = "'2.58"
$ws.Range("E41").Value = "  +3.94%  "

# Row 42
$ws.Range("E42").Value = "  +2.56%  "

# Row 43
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.74"
$ws.Range("E44").Value = "  +9.01%  "

# Row 45
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'122.52"
$ws.Range("E45").Value = "  +2.17%  "

# Row 46
$ws.Range("D46").Value = "'21.58"
$ws.Range("E46").Value = "  +1.70%  "

# Row 47
$ws.Range("E47").Value = "  +2.92%  "

# Row 48
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  +1.62%  "

# Row 49
$ws.Range("D49").Value = "2.030.11"
$ws.Range("E49").Value = "  +1.04%  "

# Row 50
$ws.Range("D50").Value = "3.314.81"
$ws.Range("E50").Value = "  +2.68%  "

# Row 51
$ws.Range("D51").Value = "'0.0319"
$ws.Range("E51").Value = "  +1.55%  "
